# Delete the old header row (row 2: "municipio"/"CASOS"/"ÓBITOS") so the
# municipal data rows shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()

# Remove the two trailing aggregate rows ("outros paises" / "outros estados")
# which, after the shift above, now live at rows 75 and 76.
$ws.Rows.Item(76).Delete()
$ws.Rows.Item(75).Delete()
